{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst notePara = lastParagraph.insertParagraph(\n  \"A note on the Hypsometric Equation used in this script. The Non-Isothermal Hypsometric Equation was collected from Wallace and Hobbs 1977 on pages 60-61. To find this equation I had to go to the Weaver Science and Engineering Library and find the original, version 1, of the book. This was in an effort to understand the associated MetPy function of \u201cheight_to_pressure_std\u201d which was originally used to calculate the pressure for the Lifting Condensation Level in my code. I have since coded the LCL pressure by hand using the equation from Wallace and Hobss instead of relying solely on MetPy.\",\n  Word.InsertLocation.after\n);\nconst asidePara = notePara.insertParagraph(\n  \" (This totally didn\u2019t take about three hours to track down\u2026 But at least it works now!)\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Grab the document's current final paragraph (the empty one right before\n# the section break) and append two new paragraphs of write-up text after it.\n$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$r = $lastPara.Range\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n\n$notePara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$notePara.Range.InsertAfter('A note on the Hypsometric Equation used in this script. The Non-Isothermal Hypsometric Equation was collected from Wallace and Hobbs 1977 on pages 60-61. To find this equation I had to go to the Weaver Science and Engineering Library and find the original, version 1, of the book. This was in an effort to understand the associated MetPy function of \u201cheight_to_pressure_std\u201d which was originally used to calculate the pressure for the Lifting Condensation Level in my code. I have since coded the LCL pressure by hand using the equation from Wallace and Hobss instead of relying solely on MetPy.')\n\n$asideParaAnchor = $d.Paragraphs.Item($d.Paragraphs.Count)\n$asideParaAnchor.Range.Collapse(0)\n$asideParaAnchor.Range.InsertParagraphAfter()\n\n$asidePara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$asidePara.Range.InsertAfter(' (This totally didn\u2019t take about three hours to track down\u2026 But at least it works now!)')\n"}
